$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.133006634347805
$ws.Range("E2").Value = 8.742833136766334
$ws.Range("F2").Value = 23.59218716732487
